# The presentation currently uses the "Integral" colour theme for its
# slide master / slides (ppt/theme/theme1.xml). The commit swaps the
# active colour theme over to the plain "Office Theme" palette that was
# already bundled in the deck (ppt/theme/theme2.xml), i.e. the classic
# Design > Themes > "Office Theme" gallery click.
#
# Re-point every themed colour slot (background/text/accent/hyperlink)
# on the presentation's colour scheme to the "Office Theme" RGB values.
# PowerPoint stores theme colours in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB() isn't available in this host, so colours are passed as the
# usual VBA-style decimal (R + G*256 + B*65536) values.

$p = $ppt.ActivePresentation

$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}

Write-Host "Applied Office Theme colour scheme to the presentation theme."
